$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that are formatted as plain text in the
# source data (e.g. "537.29", "0.0000140", thousands-separated values like
# "59.009.05", etc.). Force the whole column to Text format before writing so
# Excel does not silently reinterpret numeric-looking strings as numbers and
# normalize away trailing zeros / switch to scientific notation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '59.009.05'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '2.499.87'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '537.29'
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").Value = '137.88'
$ws.Range("E6").Value = '  -1.52%  '
$ws.Range("D8").Value = '0.566'
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("D9").Value = '2.525.06'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '0.101'
$ws.Range("E10").Value = '  +0.75%  '
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").Value = '5.36'
$ws.Range("E12").Value = '  -2.43%  '
$ws.Range("D13").Value = '0.346'
$ws.Range("E13").Value = '  -3.63%  '
$ws.Range("D14").Value = '2.954.18'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("D15").Value = '23.18'
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").Value = '58.922.42'
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").Value = '0.0000140'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").Value = '2.515.81'
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("D19").Value = '11.08'
$ws.Range("E19").Value = '  +0.68%  '
$ws.Range("D20").Value = '4.28'
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("D21").Value = '326.04'
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '5.86'
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '65.67'
$ws.Range("E24").Value = '  +5.55%  '
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").Value = '7.63'
$ws.Range("E28").Value = '  -2.35%  '
$ws.Range("D29").Value = '0.0₃0779'
$ws.Range("E29").Value = '  +1.58%  '
$ws.Range("D30").Value = '6.70'
$ws.Range("E30").Value = '  -1.99%  '
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '169.31'
$ws.Range("E32").Value = '  +4.78%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.21'
$ws.Range("E33").Value = '  +7.28%  '
$ws.Range("E34").Value = '  +2.65%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '18.57'
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("E37").Value = '  -2.89%  '
$ws.Range("D38").Value = '1.56'
$ws.Range("E38").Value = '  -1.46%  '
$ws.Range("D39").Value = '36.70'
$ws.Range("E39").Value = '  -0.67%  '
$ws.Range("D40").Value = '0.825'
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '3.64'
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '5.35'
$ws.Range("E42").Value = '  +2.06%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '285.29'
$ws.Range("E43").Value = '  +1.73%  '
$ws.Range("D44").Value = '0.995'
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '131.28'
$ws.Range("E45").Value = '  +7.64%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.605'
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("D47").Value = '10.86'
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("D48").Value = '0.0933'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").Value = '0.0513'
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("D50").Value = '0.0222'
$ws.Range("E50").Value = '  -0.58%  '
$ws.Range("D51").Value = '17.55'
$ws.Range("E51").Value = '  -0.53%  '

Write-Host "Done applying crypto price/volume updates."
